$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.676.80"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "2.645.81"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.03%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.73"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +1.80%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.83"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.13%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -2.11%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.47"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "3.110.65"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("D14").Value = "60.672.96"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "2.650.75"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("E18").Value = "  -0.59%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.97"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("E21").Value = "  +0.64%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  -0.35%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.24"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.01%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -2.09%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("E31").Value = "  +1.61%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.33"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.86%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.75"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.07%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.08"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  +5.41%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.884"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -0.33%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.85"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.17%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "306.15"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +4.49%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -0.48%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.635"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("E43").Value = "  +0.03%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.10"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0557"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  +1.14%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.86"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  +0.45%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.05"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "1.983.06"
$ws.Range("E51").Value = "  -0.85%  "
